$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3000.25
$ws.Range("I76").Value = 3000.25
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 3000.25
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = ""
$ws.Range("N76").Value = -2685.25
$ws.Range("H79").Value = 3000.25
$ws.Range("I79").Value = 3000.25
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 3000.25
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = ""
$ws.Range("N79").Value = -1908.25
$ws.Range("H137").Value = 1183364.4
$ws.Range("I137").Value = 1284.9025
$ws.Range("J137").Value = 9260907
$ws.Range("K137").Value = 3854.7075
$ws.Range("L137").Value = 27782721
$ws.Range("M137").Value = -1304.7075
$ws.Range("N137").Value = -27787821

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 27501.666
$ws.Range("I10").Value = 3000
$ws.Range("J10").Value = 39752.5
$ws.Range("K10").Value = 3000
$ws.Range("L10").Value = 39752.5
$ws.Range("M10").Value = -2830
$ws.Range("N10").Value = -40092.5
$ws.Range("H32").Value = 4266.727
$ws.Range("I32").Value = 2931.568
$ws.Range("J32").Value = 14948
$ws.Range("K32").Value = 2931.568
$ws.Range("L32").Value = 14948
$ws.Range("M32").Value = -2644.568
$ws.Range("N32").Value = -15522
$ws.Range("H88").Value = 1619.5714
$ws.Range("I88").Value = 1214.5714
$ws.Range("J88").Value = 2429.5715
$ws.Range("K88").Value = 1214.5714
$ws.Range("L88").Value = 2429.5715
$ws.Range("M88").Value = -808.5714
$ws.Range("N88").Value = -3241.5715
$ws.Range("H91").Value = 1619.5714
$ws.Range("I91").Value = 1214.5714
$ws.Range("J91").Value = 2429.5715
$ws.Range("K91").Value = 1214.5714
$ws.Range("L91").Value = 2429.5715
$ws.Range("M91").Value = 189.4286
$ws.Range("N91").Value = -5237.5715
$ws.Range("H132").Value = 3208.8667
$ws.Range("I132").Value = 2776.5715
$ws.Range("J132").Value = 3587.125
$ws.Range("K132").Value = 8329.7145
$ws.Range("L132").Value = 10761.375
$ws.Range("M132").Value = -5799.7145
$ws.Range("N132").Value = -15821.375

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1721.5652
$ws.Range("I86").Value = 1588.6666
$ws.Range("J86").Value = 2200
$ws.Range("K86").Value = 1588.6666
$ws.Range("L86").Value = 2200
$ws.Range("M86").Value = -465.6666
$ws.Range("N86").Value = -4446
$ws.Range("H89").Value = 1721.5652
$ws.Range("I89").Value = 1588.6666
$ws.Range("J89").Value = 2200
$ws.Range("K89").Value = 7943.333000000001
$ws.Range("L89").Value = 11000
$ws.Range("M89").Value = -2327.333000000001
$ws.Range("N89").Value = -22232
$ws.Range("H107").Value = 1691.6154
$ws.Range("I107").Value = 1631.1
$ws.Range("J107").Value = 1893.3334
$ws.Range("K107").Value = 1631.1
$ws.Range("L107").Value = 1893.3334
$ws.Range("M107").Value = 288.9000000000001
$ws.Range("N107").Value = -5733.3334
$ws.Range("H134").Value = 649307.4
$ws.Range("I134").Value = 956146.9
$ws.Range("J134").Value = 4944.4
$ws.Range("K134").Value = 2868440.7
$ws.Range("L134").Value = 14833.2
$ws.Range("M134").Value = -2865905.7
$ws.Range("N134").Value = -19903.2
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = ""
$ws.Range("N140").Value = 0

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 21827.455
$ws.Range("I59").Value = 4204
$ws.Range("J59").Value = 22666.666
$ws.Range("K59").Value = 4204
$ws.Range("L59").Value = 22666.666
$ws.Range("M59").Value = -3059
$ws.Range("N59").Value = -24956.666
$ws.Range("H134").Value = 1504.8959
$ws.Range("I134").Value = 1169.7297
$ws.Range("J134").Value = 2632.2727
$ws.Range("K134").Value = 3509.189100000001
$ws.Range("L134").Value = 7896.8181
$ws.Range("M134").Value = -974.1891000000005
$ws.Range("N134").Value = -12966.8181
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = ""
$ws.Range("N135").Value = 0

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 221.36363
$ws.Range("I11").Value = 119
$ws.Range("J11").Value = 400.5
$ws.Range("K11").Value = 357
$ws.Range("L11").Value = 1201.5
$ws.Range("M11").Value = -217
$ws.Range("N11").Value = -1481.5
$ws.Range("H99").Value = 3748.3333
$ws.Range("I99").Value = 3748.3333
$ws.Range("K99").Value = 11244.9999
$ws.Range("M99").Value = -8998.999899999999
$ws.Range("H114").Value = 17357814
$ws.Range("I114").Value = 50000176
$ws.Range("J114").Value = 10103957
$ws.Range("K114").Value = 150000528
$ws.Range("L114").Value = 30311871
$ws.Range("M114").Value = -149997274
$ws.Range("N114").Value = -30318379
$ws.Range("H136").Value = 3917.7856
$ws.Range("I136").Value = 2121.2856
$ws.Range("J136").Value = 5714.2856
$ws.Range("K136").Value = 6363.8568
$ws.Range("L136").Value = 17142.8568
$ws.Range("M136").Value = -1263.8568
$ws.Range("N136").Value = -27342.8568

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 1200334.2
$ws.Range("I3").Value = 1333482.5
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 1333482.5
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = -1333366.5
$ws.Range("N3").Value = -2232
$ws.Range("H70").Value = 4648.9165
$ws.Range("I70").Value = 4423.75
$ws.Range("J70").Value = 5774.75
$ws.Range("K70").Value = 4423.75
$ws.Range("L70").Value = 5774.75
$ws.Range("M70").Value = -4153.75
$ws.Range("N70").Value = -6314.75
$ws.Range("H73").Value = 4648.9165
$ws.Range("I73").Value = 4423.75
$ws.Range("J73").Value = 5774.75
$ws.Range("K73").Value = 4423.75
$ws.Range("L73").Value = 5774.75
$ws.Range("M73").Value = -3487.75
$ws.Range("N73").Value = -7646.75
$ws.Range("H80").Value = 2373.889
$ws.Range("I80").Value = 2373.5293
$ws.Range("J80").Value = 2380
$ws.Range("K80").Value = 2373.5293
$ws.Range("L80").Value = 2380
$ws.Range("M80").Value = -1375.5293
$ws.Range("N80").Value = -4376
$ws.Range("H83").Value = 2373.889
$ws.Range("I83").Value = 2373.5293
$ws.Range("J83").Value = 2380
$ws.Range("K83").Value = 11867.6465
$ws.Range("L83").Value = 11900
$ws.Range("M83").Value = -6875.646500000001
$ws.Range("N83").Value = -21884
$ws.Range("H102").Value = 1030.2646
$ws.Range("I102").Value = 942.11536
$ws.Range("J102").Value = 1316.75
$ws.Range("K102").Value = 942.11536
$ws.Range("L102").Value = 1316.75
$ws.Range("M102").Value = 679.88464
$ws.Range("N102").Value = -4560.75
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = ""
$ws.Range("N135").Value = 0
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = ""
$ws.Range("N138").Value = 0
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = ""
$ws.Range("N140").Value = 0

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3416.2131
$ws.Range("I132").Value = 2969.6956
$ws.Range("J132").Value = 4785.533
$ws.Range("K132").Value = 8909.086800000001
$ws.Range("L132").Value = 14356.599
$ws.Range("M132").Value = -6379.086800000001
$ws.Range("N132").Value = -19416.599

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 42003.5
$ws.Range("I12").Value = 14000
$ws.Range("J12").Value = 70007
$ws.Range("K12").Value = 14000
$ws.Range("L12").Value = 70007
$ws.Range("M12").Value = -13858
$ws.Range("N12").Value = -70291
$ws.Range("H132").Value = 2179.5
$ws.Range("I132").Value = 1959.3778
$ws.Range("J132").Value = 2839.8667
$ws.Range("K132").Value = 5878.1334
$ws.Range("L132").Value = 8519.6001
$ws.Range("M132").Value = -3348.1334
$ws.Range("N132").Value = -13579.6001
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = ""
$ws.Range("N139").Value = 0
$ws.Range("H140").Value = 37143.2
$ws.Range("J140").Value = 37143.2
$ws.Range("L140").Value = 37143.2
$ws.Range("N140").Value = -47503.2
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = ""
$ws.Range("N141").Value = 0
